# Case and Fatality Demographics Data Updated
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Cases by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value = 259
$ws.Range("B3").Value = 1242
$ws.Range("B4").Value = 3328
$ws.Range("B5").Value = 14267
$ws.Range("B6").Value = 15704
$ws.Range("B7").Value = 13740
$ws.Range("B8").Value = 11681
$ws.Range("B9").Value = 4221
$ws.Range("B10").Value = 2838
$ws.Range("B11").Value = 1685
$ws.Range("B12").Value = 1089
$ws.Range("B13").Value = 1700
$ws.Range("A18").Select()

# ---------------------------------------------------------------------------
# Sheet 2: Cases by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 24344
$ws.Range("B3").Value = 46529
$ws.Range("B4").Value = 895
$ws.Range("B2:B4").Select()

# ---------------------------------------------------------------------------
# Sheet 3: Cases by RaceEthnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 922
$ws.Range("B3").Value = 12018
$ws.Range("B4").Value = 27129
$ws.Range("B5").Value = 385
$ws.Range("B6").Value = 23170
$ws.Range("B7").Value = 8144
$ws.Range("A28").Select()

# ---------------------------------------------------------------------------
# Sheet 4: Fatalities by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B3").Value = 9
$ws.Range("B5").Value = 183
$ws.Range("B6").Value = 598
$ws.Range("B7").Value = 1789
$ws.Range("B8").Value = 4108
$ws.Range("B9").Value = 3459
$ws.Range("B10").Value = 4463
$ws.Range("B11").Value = 5061
$ws.Range("B12").Value = 5072
$ws.Range("B13").Value = 13356
$ws.Range("B2:B13").Select()

# ---------------------------------------------------------------------------
# Sheet 5: Fatalities by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 16072
$ws.Range("B3").Value = 22055

# ---------------------------------------------------------------------------
# Sheet 6: Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 735
$ws.Range("B3").Value = 3629
$ws.Range("B4").Value = 17824
$ws.Range("B5").Value = 200
$ws.Range("B6").Value = 15719
$ws.Range("C22").Select()

# Re-select the first sheet as tab-selected / active like the original file
$wb.Worksheets.Item("Cases by Age Group").Activate()
